$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 168
$ws.Range("F6").Value = 4945
$ws.Range("F7").Value = 385
$ws.Range("F8").Value = 565
$ws.Range("F9").Value = 864
$ws.Range("F11").Value = 63
$ws.Range("F13").Value = 538
$ws.Range("F16").Value = 1611
$ws.Range("F17").Value = 1423
$ws.Range("F18").Value = 703
$ws.Range("F20").Value = 172
$ws.Range("F21").Value = 257
$ws.Range("F22").Value = 465
$ws.Range("F23").Value = 109
$ws.Range("F24").Value = 1036
$ws.Range("F27").Value = 1674
$ws.Range("F28").Value = 144
$ws.Range("F29").Value = 78
$ws.Range("F30").Value = 12
$ws.Range("F31").Value = 194
$ws.Range("F34").Value = 7
$ws.Range("F37").Value = 548
$ws.Range("F38").Value = 71
$ws.Range("F39").Value = 15
$ws.Range("F40").Value = 18
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 133
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 168
$ws.Range("F8").Value = 4945
$ws.Range("F9").Value = 385
$ws.Range("F10").Value = 565
$ws.Range("F12").Value = 133
$ws.Range("F13").Value = 864
$ws.Range("F17").Value = 63
$ws.Range("F19").Value = 538
$ws.Range("F23").Value = 1611
$ws.Range("F24").Value = 1423
$ws.Range("F25").Value = 703
$ws.Range("F27").Value = 172
$ws.Range("F28").Value = 257
$ws.Range("F30").Value = 465
$ws.Range("F31").Value = 109
$ws.Range("F32").Value = 1036
$ws.Range("F34").Value = 1674
$ws.Range("F35").Value = 144
$ws.Range("F36").Value = 78
$ws.Range("F37").Value = 12
$ws.Range("F38").Value = 194
$ws.Range("F41").Value = 7
$ws.Range("F43").Value = 548
$ws.Range("F44").Value = 71
$ws.Range("F45").Value = 15
$ws.Range("F46").Value = 18
